# Update xlsx datetime test data
#
# The "datetime" sheet's sample row (row 2) held timestamps formatted with a
# -07:00 / MST offset. The refreshed fixture was generated on a UTC host, so
# those same logical values now render with a Z / +0000 / UTC offset
# instead. Row 1 (the format-name headers) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datetime")

$ws.Range("D2").Value = "1989-11-09T15:17:59.123-00:00"
$ws.Range("F2").Value = "Thu, 09 Nov 1989 15:17:59 UTC"
$ws.Range("G2").Value = "Thu, 09 Nov 1989 15:17:59 +0000"
$ws.Range("H2").Value = "1989-11-09T15:17:59Z"
$ws.Range("I2").Value = "1989-11-09T15:17:59.1234567Z"
$ws.Range("K2").Value = "1989-11-09T15:17:59Z"
$ws.Range("L2").Value = "09 Nov 89 15:17 UTC"
$ws.Range("M2").Value = "09 Nov 89 15:17 +0000"
$ws.Range("N2").Value = "Thursday, 09-Nov-89 15:17:59 UTC"
$ws.Range("O2").Value = "Thu Nov 09 15:17:59 +0000 1989"
$ws.Range("T2").Value = "Thu Nov  9 15:17:59 UTC 1989"

# Selection moved from D2 to T2, with the view scrolled so column K is
# leftmost (as far as that's representable through the object model).
$ws.Range("T2").Select()
